$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 367
$ws.Range("C3").Value = 367
$ws.Range("C4").Value = 367
$ws.Range("C5").Value = 367
$ws.Range("C6").Value = 367
$ws.Range("C7").Value = 872
$ws.Range("C8").Value = 899
$ws.Range("C9").Value = 905
$ws.Range("C10").Value = 892
$ws.Range("C11").Value = 907
$ws.Range("C12").Value = 1102
$ws.Range("C13").Value = 749
$ws.Range("C14").Value = 676
$ws.Range("C15").Value = 681
